$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the report title (October 2016 -> November 2016)
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Table 4.3. Receipts, Average Cost, and Quality of Fossil Fuels: Independent Power Producers, 2006 - November 2016 (continued)"

# ------------------------------------------------------------------
# 2. Insert a new "November" data row for the Natural Gas - Year 2016
#    block. This pushes the existing "Year to Date" section (and
#    everything below it) down by one row.
# ------------------------------------------------------------------
$ws.Rows.Item(53).Insert()

# Copy the formatting of the row above (October's data row, now row 52)
# down onto the newly inserted blank row so the new row matches the
# other month rows exactly.
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new November row (Natural Gas, Year 2016)
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 1294
$ws.Range("C53").Value = 46
$ws.Range("D53").Value = "W"
$ws.Range("E53").Value = "W"
$ws.Range("F53").Value = 5.43
$ws.Range("G53").Value = 83.4
$ws.Range("H53").Value = 317208
$ws.Range("I53").Value = 307453
$ws.Range("J53").Value = 2.6
$ws.Range("K53").Value = 2.68
$ws.Range("L53").Value = 93.3
$ws.Range("M53").Value = "W"

# ------------------------------------------------------------------
# 3. Refresh the "Year to Date" figures (now rows 55-57) to include
#    the extra month of data.
# ------------------------------------------------------------------
# Year to Date - 2014
$ws.Range("B55").Value = 12261
$ws.Range("C55").Value = 435
$ws.Range("D55").Value = 2.48
$ws.Range("E55").Value = 70.34
$ws.Range("F55").Value = 5.36
$ws.Range("G55").Value = 70.2
$ws.Range("H55").Value = 3732483
$ws.Range("I55").Value = 3623155
$ws.Range("J55").Value = 4.97
$ws.Range("K55").Value = 5.12
$ws.Range("L55").Value = 92.6

# Year to Date - 2015
$ws.Range("B56").Value = 12809
$ws.Range("C56").Value = 458
$ws.Range("D56").Value = 2.45
$ws.Range("E56").Value = 68.46
$ws.Range("F56").Value = 5.2
$ws.Range("G56").Value = 61.7
$ws.Range("H56").Value = 4320981
$ws.Range("I56").Value = 4180105
$ws.Range("J56").Value = 3
$ws.Range("K56").Value = 3.1
$ws.Range("L56").Value = 93.2

# Year to Date - 2016
$ws.Range("B57").Value = 12066
$ws.Range("C57").Value = 437
$ws.Range("D57").Value = 2.5
$ws.Range("E57").Value = 69.01
$ws.Range("F57").Value = 5.43
$ws.Range("G57").Value = 68.4
$ws.Range("H57").Value = 4462966
$ws.Range("I57").Value = 4316582
$ws.Range("J57").Value = 2.45
$ws.Range("K57").Value = 2.54
$ws.Range("L57").Value = 93.7

# ------------------------------------------------------------------
# 4. Rename the "Rolling 12 Months Ending in October" header (now row
#    58) to "...November" and refresh its figures (rows 59-60).
# ------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# Rolling 12 months - 2015
$ws.Range("B59").Value = 14329
$ws.Range("C59").Value = 511
$ws.Range("F59").Value = 5.19
$ws.Range("G59").Value = 63
$ws.Range("H59").Value = 4643038
$ws.Range("I59").Value = 4491621
$ws.Range("J59").Value = 3.07
$ws.Range("K59").Value = 3.18

# Rolling 12 months - 2016
$ws.Range("B60").Value = 13807
$ws.Range("C60").Value = 503
$ws.Range("F60").Value = 5.47
$ws.Range("G60").Value = 74.4
$ws.Range("H60").Value = 4825275
$ws.Range("I60").Value = 4666672
$ws.Range("J60").Value = 2.44
$ws.Range("K60").Value = 2.52

Write-Output "Applied November 2016 update"
